# Generate Report for Handback
# -----------------------------------------------------------------------
# This script reproduces (logically) the OOXML diff:
#  1. Overview sheet status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (a shared string used by the
#     Status columns on every sheet -> updating every cell that held the
#     old text collapses back onto one shared string).
#  2. zh-cn sheet: rows 2 & 3 get their "Latest Target File" (I) and
#     "Latest Handback File" (J) populated (I becomes a hyperlink to the
#     source .md file, like column A already is), and the shared
#     "Latest Handback DateTime" (K) flips from the epoch placeholder to
#     a real timestamp.
#  3. de-de sheet: same shape of edit as zh-cn, but with its own
#     xlf filenames / timestamp.
#  4. Column widths on the affected columns widen to fit the new text.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$urlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/192bcff187947b0e07e5e386917f4000502f5300/e2e"
$mdFile1 = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md"
$mdFile2 = "ba774427-4f78-4031-ad1a-bc070f21edd8.md"

# ---------------------------------------------------------------------
# 1. Overview sheet - Status text for both rows/both language columns.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Widen columns E (zh-cn status) and F (de-de status) to fit the longer text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2 (90d90118... file): same Status text as above, propagated too.
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Latest Target File (I) + hyperlink, Latest Handback File (J)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "$urlBase/$mdFile1", "", "", $mdFile1) | Out-Null
$wsZh.Range("J2").Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.246882f00340f1c95e596140032e420920c74481.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "$urlBase/$mdFile2", "", "", $mdFile2) | Out-Null
$wsZh.Range("J3").Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.1b32bcab9903cb6e6644e907e97d1bfe70dbd639.zh-cn.xlf"

# Latest Handback DateTime (K) - shared by both rows; was "0001-01-01 00:00:00".
$wsZh.Range("K2").Value = "2016-09-07 08:08:41"
$wsZh.Range("K3").Value = "2016-09-07 08:08:41"

# Widen Status (C) and the two newly-filled columns (I, J).
$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "$urlBase/$mdFile1", "", "", $mdFile1) | Out-Null
$wsDe.Range("J2").Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.246882f00340f1c95e596140032e420920c74481.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "$urlBase/$mdFile2", "", "", $mdFile2) | Out-Null
$wsDe.Range("J3").Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.1b32bcab9903cb6e6644e907e97d1bfe70dbd639.de-de.xlf"

# Latest Handback DateTime (K) - shared by both rows; was "0001-01-01 00:00:00".
$wsDe.Range("K2").Value = "2016-09-07 08:08:59"
$wsDe.Range("K3").Value = "2016-09-07 08:08:59"

$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
